$wb = $excel.ActiveWorkbook

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")

# Source citation block (B3:B7) gets replaced with the new reference details
$about.Range("B3").Value = "Massachusetts Institute of Technology"
$about.Range("B4").Value = 2021
$about.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$about.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$about.Range("B7").Value = "Abstract"

# Old note in C8 is cleared (no longer applicable)
$about.Range("C8").ClearContents()

# New note added below, in A9
$about.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# Remove the old chart image that illustrated the previous (BNEF) source
if ($about.Shapes.Count -gt 0) {
    for ($i = $about.Shapes.Count; $i -ge 1; $i--) {
        $about.Shapes.Item($i).Delete()
    }
}

# --- "PDiBCpDoC" sheet updates ---
$data = $wb.Worksheets.Item("PDiBCpDoC")

# B2 now derives from the new source's quoted learning-rate range instead of a hardcoded constant
$data.Range("B2").Formula = "=AVERAGE(0.2,0.27)"
